$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.614100575447083
$ws.Range("B1").Value = 2.748774290084839
$ws.Range("C1").Value = 3.178627967834473
$ws.Range("D1").Value = 3.556381464004517
$ws.Range("E1").Value = 1.619780302047729
